$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column before column A, shifting existing data (A:E) to (B:F)
$ws.Range("A:A").Insert()

# Rename the shifted header cells (B1:F1) to their new labels
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Set the new header in A1 and the method names in A2:A9
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Adjust column widths to match the target layout (best-fit widths)
$ws.Range("A:A").ColumnWidth = 12.451822916666666
$ws.Range("B:B").ColumnWidth = 3.1666666666666665
$ws.Range("C:C").ColumnWidth = 2.3072916666666665
$ws.Range("D:F").ColumnWidth = 11.166666666666666
